$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: Title "A Table, with a caption" -- consolidate word+space runs.
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 2).Text = "A "
$titleRange.Characters(3, 7).Text = "Table, "
$titleRange.Characters(10, 5).Text = "with "
$titleRange.Characters(15, 2).Text = "a "
# "caption" (chars 17-23) is already its own run; leave untouched.

# Shape 3: TextBox "Demonstration of simple table syntax, with alignment"
# -- consolidate word+space runs.
$captionShape = $s.Shapes.Item(3)
$captionRange = $captionShape.TextFrame.TextRange
$captionRange.Characters(1, 14).Text = "Demonstration "
$captionRange.Characters(15, 3).Text = "of "
$captionRange.Characters(18, 7).Text = "simple "
$captionRange.Characters(25, 6).Text = "table "
$captionRange.Characters(31, 8).Text = "syntax, "
$captionRange.Characters(39, 5).Text = "with "
# "alignment" (chars 44-52) is already its own run; leave untouched.
